$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: update timestart/timeend and fill in MOVIL / phone1
$ws.Range("F2").Value = 1725958800
$ws.Range("G2").Value = 1728118800

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "3209006290"
$ws.Range("H2").Style = "Normal"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "573209006290"
$ws.Range("I2").Style = "Normal"

# Row 3: update timestart/timeend
$ws.Range("F3").Value = 1725958800
$ws.Range("G3").Value = 1725958800

# Row 4: update timestart/timeend
$ws.Range("F4").Value = 1725958800
$ws.Range("G4").Value = 1726563600
